$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.260.29'
$ws.Range("E2").Value = '  +0.63%  '
$ws.Range("D3").Value = '2.443.85'
$ws.Range("E3").Value = '  -0.08%  '
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '0.998'
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.28%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '571.93'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.79%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '146.89'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +0.39%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +0.03%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.540'
$c.Style = "Normal"
$ws.Range("E8").Value = '  +0.99%  '
$ws.Range("D9").Value = '2.438.00'
$ws.Range("E9").Value = '  -0.50%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.111'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -0.80%  '
$ws.Range("E11").Value = '  +1.05%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '5.24'
$c.Style = "Normal"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.355'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -0.08%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '27.14'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +0.27%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.0000179'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -2.34%  '
$ws.Range("D16").Value = '2.876.28'
$ws.Range("E16").Value = '  -0.43%  '
$ws.Range("D17").Value = '62.971.11'
$ws.Range("E17").Value = '  +0.56%  '
$ws.Range("D18").Value = '2.447.79'
$ws.Range("E18").Value = '  -0.02%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '11.33'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +0.22%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '7.37'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +5.86%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '327.83'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +1.10%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '4.19'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +0.23%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '2.08'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +12.81%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '1.04'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +4.46%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '65.24'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -3.29%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '619.94'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +5.29%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '8.86'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +1.56%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '0.0000103'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +1.19%  '
$ws.Range("B29").Value = 'Fetch.AI'
$ws.Range("C29").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '1.51'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +3.67%  '
$ws.Range("B30").Value = 'WrappedeETH'
$ws.Range("C30").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D30").Value = '2.558.33'
$ws.Range("E30").Value = '  -0.20%  '
$ws.Range("B31").Value = 'Binance-PegBSC-USD'
$ws.Range("C31").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +0.32%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '8.29'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -2.29%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.142'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -3.95%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '1.90'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +1.04%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '5.21'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +6.88%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '1.53'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -1.82%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +0.04%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.381'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -0.72%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '5.44'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -0.44%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '18.74'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -0.52%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '146.59'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -0.99%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '2.70'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +10.86%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '1.80'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -1.90%  '
$ws.Range("E44").Value = '  -0.51%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '149.04'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -0.40%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '3.75'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +1.93%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '21.36'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +3.60%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.0536'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -0.19%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.598'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -0.89%  '
$ws.Range("E50").Value = '  +0.33%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.0916'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -1.14%  '
